$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<quebec>"

# Row 3
$ws.Range("B3").Value = "<an>"
$ws.Range("C3").Value = 53

# Row 4
$ws.Range("B4").Value = "<mike>"
$ws.Range("C4").Value = 52

# Row 5
$ws.Range("C5").Value = 51

# Row 6
$ws.Range("C6").Value = 51

# Row 7
$ws.Range("C7").Value = 49

# Row 8
$ws.Range("B8").Value = "<they>"
$ws.Range("C8").Value = 49

# Row 9
$ws.Range("B9").Value = "<yes>"
$ws.Range("C9").Value = 48

# Row 11
$ws.Range("C11").Value = 48

# Row 12
$ws.Range("B12").Value = "<downward>"
$ws.Range("C12").Value = 50

# Row 13
$ws.Range("C13").Value = 45

# Row 14
$ws.Range("B14").Value = "<my>"
$ws.Range("C14").Value = 44

# Row 16
$ws.Range("B16").Value = "<like>"
$ws.Range("C16").Value = 50

# Row 17
$ws.Range("C17").Value = 58

# Row 18
$ws.Range("C18").Value = 43
